{"js": "// Fill in the empty AdmNo/Name row of the group-members table with the\n// new group member's admin number and name (bold, matching the other rows).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Row 0 is the header (\"AdmNo\" / \"Name\"); row 1 is the first (empty) data\n// row that needs the new member's details.\nconst table = tables.items[0];\nconst admNoCell = table.getCell(1, 0);\nconst nameCell = table.getCell(1, 1);\n\n// Insert the text into the cell's existing (empty) paragraph rather than\n// replacing the whole cell body, so the paragraph itself is preserved.\nconst admNoPara = admNoCell.body.paragraphs.getFirst();\nconst namePara = nameCell.body.paragraphs.getFirst();\nadmNoPara.insertText(\"2112617\", \"Replace\");\nnamePara.insertText(\"Lim Ke Zhen Joaquin\", \"Replace\");\n\n// The cell paragraphs already carry a bold paragraph mark, but match the\n// other populated rows by making the inserted run itself bold too.\nadmNoPara.font.bold = true;\nnamePara.font.bold = true;\n\nawait context.sync();\n", "ps1": "# Fill in the empty AdmNo/Name row of the group-members table with the\n# new group member's admin number and name (bold, matching the other rows).\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Row 1 is the header (\"AdmNo\" / \"Name\"); row 2 is the first (empty) data\n# row that needs the new member's details. Table.Cell is 1-based.\n$admNoCell = $table.Cell(2, 1)\n$nameCell = $table.Cell(2, 2)\n\n$admNoCell.Range.Text = \"2112617\"\n$nameCell.Range.Text = \"Lim Ke Zhen Joaquin\"\n\n# The cell paragraphs already carry a bold paragraph mark, but match the\n# other populated rows by making the inserted run itself bold too.\n$admNoCell.Range.Font.Bold = 1\n$nameCell.Range.Font.Bold = 1\n"}
